$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 4.970700000000004
$ws.Range("C4").Value = -14.1287

$ws.Range("C5").Value = -14.78700000000001

$ws.Range("B7").Value = 6.189399999999997

$ws.Range("C8").Value = -11.75539999999999

$ws.Range("B16").Value = 9.132300000000008
$ws.Range("C16").Value = -12.08100000000001
